# "Initial commit of Component Required Field Validation"
#
# Renames the worksheet, rewrites the header row (row 1) with the new set
# of required-field column headers (inserted ModelYear/ManufacturerID/
# ModelID earlier in the row, fixed the AssetCategoryID -> AssetcategoryID
# casing, moved PreferredPMShift ahead of StationLocation, and appended a
# new WorkOrders column), applies Text formatting + Bold to the header
# row, autosizes the columns, and sets the page orientation to portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: "Reqd Field Validation" -> "RequiredFields"
$ws.Name = "RequiredFields"

# Final header order (A1:X1)
$headers = @(
  "AssetID","Source","LegacyIDSource","LegacyID","Description",
  "ModelYear","ManufacturerID","ModelID",
  "EquipmentType","PMProgramType",
  "MeterTypesClass","Maintenance","PMClass","Standards",
  "RentalRates","Resources","AssetcategoryID","AssignedPM","AssignedRepair",
  "PreferredPMShift","StationLocation",
  "DepartmentID","DepartmentForPM","WorkOrders"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$lastCol = $headers.Length
$headerRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(1, $lastCol))

# Header formatting: Text number format + Bold font
$headerRange.NumberFormat = "@"
$headerRange.Font.Bold = $true

# Autofit columns to the new header content
$headerRange.EntireColumn.AutoFit() | Out-Null

# Page layout: portrait orientation
$ws.PageSetup.Orientation = 1
